$d = $word.ActiveDocument

$d.Content.Find.Execute("70-48=22", $true, $false, $false, $false, $false, $true, 1, $false, "42-32=10", 2) | Out-Null
$d.Content.Find.Execute("60-24=36", $true, $false, $false, $false, $false, $true, 1, $false, "6+62=68", 2) | Out-Null
$d.Content.Find.Execute("47-5=42", $true, $false, $false, $false, $false, $true, 1, $false, "20+40=60", 2) | Out-Null
$d.Content.Find.Execute("14+24=38", $true, $false, $false, $false, $false, $true, 1, $false, "25-16=9", 2) | Out-Null
$d.Content.Find.Execute("66-39=27", $true, $false, $false, $false, $false, $true, 1, $false, "25+10=35", 2) | Out-Null
$d.Content.Find.Execute("37+35=72", $true, $false, $false, $false, $false, $true, 1, $false, "12-1=11", 2) | Out-Null
$d.Content.Find.Execute("66+2=68", $true, $false, $false, $false, $false, $true, 1, $false, "68-51=17", 2) | Out-Null
$d.Content.Find.Execute("22+5=27", $true, $false, $false, $false, $false, $true, 1, $false, "63+20=83", 2) | Out-Null
$d.Content.Find.Execute("88-62=26", $true, $false, $false, $false, $false, $true, 1, $false, "38-33=5", 2) | Out-Null
$d.Content.Find.Execute("36+13=49", $true, $false, $false, $false, $false, $true, 1, $false, "96-39=57", 2) | Out-Null
$d.Content.Find.Execute("67-9=58", $true, $false, $false, $false, $false, $true, 1, $false, "36-9=27", 2) | Out-Null
$d.Content.Find.Execute("60-2=58", $true, $false, $false, $false, $false, $true, 1, $false, "20+56=76", 2) | Out-Null
$d.Content.Find.Execute("60-25=35", $true, $false, $false, $false, $false, $true, 1, $false, "9+58=67", 2) | Out-Null
$d.Content.Find.Execute("52-42=10", $true, $false, $false, $false, $false, $true, 1, $false, "86-65=21", 2) | Out-Null
$d.Content.Find.Execute("31-29=2", $true, $false, $false, $false, $false, $true, 1, $false, "24+54=78", 2) | Out-Null
$d.Content.Find.Execute("6+69=75", $true, $false, $false, $false, $false, $true, 1, $false, "64+20=84", 2) | Out-Null
$d.Content.Find.Execute("57-15=42", $true, $false, $false, $false, $false, $true, 1, $false, "8+36=44", 2) | Out-Null
$d.Content.Find.Execute("89-3=86", $true, $false, $false, $false, $false, $true, 1, $false, "84-27=57", 2) | Out-Null
$d.Content.Find.Execute("41+0=41", $true, $false, $false, $false, $false, $true, 1, $false, "94-10=84", 2) | Out-Null
$d.Content.Find.Execute("17+5=22", $true, $false, $false, $false, $false, $true, 1, $false, "68-45=23", 2) | Out-Null
$d.Content.Find.Execute("5+22=27", $true, $false, $false, $false, $false, $true, 1, $false, "60-53=7", 2) | Out-Null
$d.Content.Find.Execute("24+15=39", $true, $false, $false, $false, $false, $true, 1, $false, "62+1=63", 2) | Out-Null
$d.Content.Find.Execute("79+16=95", $true, $false, $false, $false, $false, $true, 1, $false, "7+85=92", 2) | Out-Null
$d.Content.Find.Execute("50-0=50", $true, $false, $false, $false, $false, $true, 1, $false, "71-66=5", 2) | Out-Null
$d.Content.Find.Execute("84-21=63", $true, $false, $false, $false, $false, $true, 1, $false, "53+0=53", 2) | Out-Null
$d.Content.Find.Execute("96-75=21", $true, $false, $false, $false, $false, $true, 1, $false, "72-32=40", 2) | Out-Null
$d.Content.Find.Execute("48+50=98", $true, $false, $false, $false, $false, $true, 1, $false, "92-60=32", 2) | Out-Null
$d.Content.Find.Execute("88-26=62", $true, $false, $false, $false, $false, $true, 1, $false, "44-17=27", 2) | Out-Null
$d.Content.Find.Execute("67-31=36", $true, $false, $false, $false, $false, $true, 1, $false, "32+0=32", 2) | Out-Null
$d.Content.Find.Execute("29+36=65", $true, $false, $false, $false, $false, $true, 1, $false, "26-16=10", 2) | Out-Null
$d.Content.Find.Execute("85-50=35", $true, $false, $false, $false, $false, $true, 1, $false, "22+24=46", 2) | Out-Null
$d.Content.Find.Execute("50+42=92", $true, $false, $false, $false, $false, $true, 1, $false, "80-60=20", 2) | Out-Null
$d.Content.Find.Execute("35+18=53", $true, $false, $false, $false, $false, $true, 1, $false, "44+20=64", 2) | Out-Null
$d.Content.Find.Execute("55+35=90", $true, $false, $false, $false, $false, $true, 1, $false, "32+22=54", 2) | Out-Null
$d.Content.Find.Execute("23+28=51", $true, $false, $false, $false, $false, $true, 1, $false, "49-19=30", 2) | Out-Null
$d.Content.Find.Execute("29+53=82", $true, $false, $false, $false, $false, $true, 1, $false, "54+44=98", 2) | Out-Null
$d.Content.Find.Execute("26+31=57", $true, $false, $false, $false, $false, $true, 1, $false, "66+19=85", 2) | Out-Null
$d.Content.Find.Execute("47-6=41", $true, $false, $false, $false, $false, $true, 1, $false, "28+37=65", 2) | Out-Null
$d.Content.Find.Execute("1+78=79", $true, $false, $false, $false, $false, $true, 1, $false, "17-17=0", 2) | Out-Null
$d.Content.Find.Execute("14+44=58", $true, $false, $false, $false, $false, $true, 1, $false, "89+7=96", 2) | Out-Null
$d.Content.Find.Execute("28+55=83", $true, $false, $false, $false, $false, $true, 1, $false, "70-21=49", 2) | Out-Null
$d.Content.Find.Execute("16+81=97", $true, $false, $false, $false, $false, $true, 1, $false, "67-40=27", 2) | Out-Null
$d.Content.Find.Execute("30+68=98", $true, $false, $false, $false, $false, $true, 1, $false, "11+87=98", 2) | Out-Null
$d.Content.Find.Execute("98-62=36", $true, $false, $false, $false, $false, $true, 1, $false, "60+0=60", 2) | Out-Null
$d.Content.Find.Execute("42-11=31", $true, $false, $false, $false, $false, $true, 1, $false, "32+17=49", 2) | Out-Null
$d.Content.Find.Execute("80+19=99", $true, $false, $false, $false, $false, $true, 1, $false, "59+5=64", 2) | Out-Null
$d.Content.Find.Execute("60-52=8", $true, $false, $false, $false, $false, $true, 1, $false, "87-8=79", 2) | Out-Null
$d.Content.Find.Execute("73-20=53", $true, $false, $false, $false, $false, $true, 1, $false, "5+30=35", 2) | Out-Null
$d.Content.Find.Execute("31+3=34", $true, $false, $false, $false, $false, $true, 1, $false, "36+47=83", 2) | Out-Null
$d.Content.Find.Execute("4+83=87", $true, $false, $false, $false, $false, $true, 1, $false, "82-40=42", 2) | Out-Null
$d.Content.Find.Execute("43+5=48", $true, $false, $false, $false, $false, $true, 1, $false, "11+75=86", 2) | Out-Null
$d.Content.Find.Execute("67-37=30", $true, $false, $false, $false, $false, $true, 1, $false, "75+22=97", 2) | Out-Null
$d.Content.Find.Execute("81+8=89", $true, $false, $false, $false, $false, $true, 1, $false, "43-13=30", 2) | Out-Null
$d.Content.Find.Execute("10+19=29", $true, $false, $false, $false, $false, $true, 1, $false, "1+41=42", 2) | Out-Null
$d.Content.Find.Execute("76-62=14", $true, $false, $false, $false, $false, $true, 1, $false, "29+11=40", 2) | Out-Null
$d.Content.Find.Execute("80-71=9", $true, $false, $false, $false, $false, $true, 1, $false, "81-28=53", 2) | Out-Null
$d.Content.Find.Execute("98-48=50", $true, $false, $false, $false, $false, $true, 1, $false, "78-7=71", 2) | Out-Null
$d.Content.Find.Execute("13+65=78", $true, $false, $false, $false, $false, $true, 1, $false, "63-42=21", 2) | Out-Null
$d.Content.Find.Execute("36+22=58", $true, $false, $false, $false, $false, $true, 1, $false, "62-45=17", 2) | Out-Null
$d.Content.Find.Execute("57+40=97", $true, $false, $false, $false, $false, $true, 1, $false, "54-37=17", 2) | Out-Null
$d.Content.Find.Execute("16+52=68", $true, $false, $false, $false, $false, $true, 1, $false, "71+7=78", 2) | Out-Null
$d.Content.Find.Execute("81-42=39", $true, $false, $false, $false, $false, $true, 1, $false, "64+8=72", 2) | Out-Null
$d.Content.Find.Execute("76-71=5", $true, $false, $false, $false, $false, $true, 1, $false, "4+43=47", 2) | Out-Null
$d.Content.Find.Execute("15+18=33", $true, $false, $false, $false, $false, $true, 1, $false, "82-56=26", 2) | Out-Null
$d.Content.Find.Execute("74-43=31", $true, $false, $false, $false, $false, $true, 1, $false, "43-34=9", 2) | Out-Null
$d.Content.Find.Execute("81-39=42", $true, $false, $false, $false, $false, $true, 1, $false, "74+15=89", 2) | Out-Null
$d.Content.Find.Execute("51+4=55", $true, $false, $false, $false, $false, $true, 1, $false, "84+14=98", 2) | Out-Null
$d.Content.Find.Execute("79+9=88", $true, $false, $false, $false, $false, $true, 1, $false, "16+28=44", 2) | Out-Null
$d.Content.Find.Execute("18+49=67", $true, $false, $false, $false, $false, $true, 1, $false, "3+3=6", 2) | Out-Null
$d.Content.Find.Execute("84+1=85", $true, $false, $false, $false, $false, $true, 1, $false, "80+4=84", 2) | Out-Null
$d.Content.Find.Execute("69-56=13", $true, $false, $false, $false, $false, $true, 1, $false, "56-6=50", 2) | Out-Null
$d.Content.Find.Execute("61-18=43", $true, $false, $false, $false, $false, $true, 1, $false, "51+31=82", 2) | Out-Null
$d.Content.Find.Execute("48+27=75", $true, $false, $false, $false, $false, $true, 1, $false, "83-54=29", 2) | Out-Null
$d.Content.Find.Execute("65+13=78", $true, $false, $false, $false, $false, $true, 1, $false, "33-9=24", 2) | Out-Null
$d.Content.Find.Execute("27+61=88", $true, $false, $false, $false, $false, $true, 1, $false, "86-52=34", 2) | Out-Null
$d.Content.Find.Execute("25-23=2", $true, $false, $false, $false, $false, $true, 1, $false, "51+27=78", 2) | Out-Null
$d.Content.Find.Execute("62+34=96", $true, $false, $false, $false, $false, $true, 1, $false, "81-30=51", 2) | Out-Null
$d.Content.Find.Execute("84-62=22", $true, $false, $false, $false, $false, $true, 1, $false, "38+48=86", 2) | Out-Null
$d.Content.Find.Execute("96-34=62", $true, $false, $false, $false, $false, $true, 1, $false, "87-77=10", 2) | Out-Null
$d.Content.Find.Execute("87-79=8", $true, $false, $false, $false, $false, $true, 1, $false, "8+37=45", 2) | Out-Null
$d.Content.Find.Execute("5+26=31", $true, $false, $false, $false, $false, $true, 1, $false, "97-61=36", 2) | Out-Null
$d.Content.Find.Execute("49+10=59", $true, $false, $false, $false, $false, $true, 1, $false, "25-22=3", 2) | Out-Null
$d.Content.Find.Execute("79-20=59", $true, $false, $false, $false, $false, $true, 1, $false, "39+28=67", 2) | Out-Null
$d.Content.Find.Execute("45-11=34", $true, $false, $false, $false, $false, $true, 1, $false, "70-10=60", 2) | Out-Null
$d.Content.Find.Execute("15+64=79", $true, $false, $false, $false, $false, $true, 1, $false, "96+0=96", 2) | Out-Null
$d.Content.Find.Execute("45-37=8", $true, $false, $false, $false, $false, $true, 1, $false, "52+13=65", 2) | Out-Null
$d.Content.Find.Execute("17+21=38", $true, $false, $false, $false, $false, $true, 1, $false, "93-87=6", 2) | Out-Null
$d.Content.Find.Execute("80-2=78", $true, $false, $false, $false, $false, $true, 1, $false, "15+34=49", 2) | Out-Null
$d.Content.Find.Execute("47+39=86", $true, $false, $false, $false, $false, $true, 1, $false, "39-20=19", 2) | Out-Null
$d.Content.Find.Execute("94+1=95", $true, $false, $false, $false, $false, $true, 1, $false, "78-2=76", 2) | Out-Null
$d.Content.Find.Execute("90-5=85", $true, $false, $false, $false, $false, $true, 1, $false, "19+25=44", 2) | Out-Null
$d.Content.Find.Execute("60+25=85", $true, $false, $false, $false, $false, $true, 1, $false, "74-49=25", 2) | Out-Null
$d.Content.Find.Execute("73-50=23", $true, $false, $false, $false, $false, $true, 1, $false, "84-68=16", 2) | Out-Null
$d.Content.Find.Execute("59-8=51", $true, $false, $false, $false, $false, $true, 1, $false, "55-38=17", 2) | Out-Null
$d.Content.Find.Execute("17+46=63", $true, $false, $false, $false, $false, $true, 1, $false, "22+6=28", 2) | Out-Null
$d.Content.Find.Execute("9+15=24", $true, $false, $false, $false, $false, $true, 1, $false, "83-67=16", 2) | Out-Null
$d.Content.Find.Execute("72-57=15", $true, $false, $false, $false, $false, $true, 1, $false, "77-72=5", 2) | Out-Null
$d.Content.Find.Execute("94+5=99", $true, $false, $false, $false, $false, $true, 1, $false, "62+7=69", 2) | Out-Null
$d.Content.Find.Execute("15+25=40", $true, $false, $false, $false, $false, $true, 1, $false, "38+34=72", 2) | Out-Null
$d.Content.Find.Execute("7+92=99", $true, $false, $false, $false, $false, $true, 1, $false, "2+26=28", 2) | Out-Null
